$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update changed values in rows 2-9 (columns B..N) ---

# Row 2
$ws.Range("B2").Value = 126164.3856688309
$ws.Range("C2").Value = 173691.2018284847
$ws.Range("D2").Value = 50888.46396309303
$ws.Range("E2").Value = 106789.5218456779
$ws.Range("F2").Value = 32647.55817709521
$ws.Range("G2").Value = 63410.96103557143
$ws.Range("L2").Value = 107832.1031784475
$ws.Range("M2").Value = 482158
$ws.Range("N2").Value = -1771974.217653017

# Row 3
$ws.Range("B3").Value = 370356.4476432937
$ws.Range("C3").Value = 1123642.666942821
$ws.Range("D3").Value = 222694.2805031927
$ws.Range("E3").Value = 542070.9584863047
$ws.Range("F3").Value = 407309.9862352163
$ws.Range("G3").Value = 64568.89510669051
$ws.Range("L3").Value = 539527.9978741551
$ws.Range("M3").Value = 1908521
$ws.Range("N3").Value = -3376968.119688379

# Row 4
$ws.Range("B4").Value = 202321.944253508
$ws.Range("C4").Value = 1032789.818620056
$ws.Range("D4").Value = 144963.7364611912
$ws.Range("E4").Value = 497911.3355670046
$ws.Range("F4").Value = 359063.1260869047
$ws.Range("G4").Value = 17093.31363413405
$ws.Range("L4").Value = 432653.1324859021
$ws.Range("M4").Value = 1889497
$ws.Range("N4").Value = -3598328.473917802

# Row 5
$ws.Range("B5").Value = 1319219.543874883
$ws.Range("C5").Value = 9432529.052794624
$ws.Range("D5").Value = 578261.8886768896
$ws.Range("E5").Value = 5022680.103990881
$ws.Range("F5").Value = 1225877.166913188
$ws.Range("G5").Value = 2189371.891411028
$ws.Range("L5").Value = 3835450.592618194
$ws.Range("M5").Value = 16592698
$ws.Range("N5").Value = -20436089.77884024

# Row 6
$ws.Range("B6").Value = 2678802.235000358
$ws.Range("C6").Value = 11985090.93287172
$ws.Range("D6").Value = 120216.5409299204
$ws.Range("E6").Value = 1899651.385322116
$ws.Range("F6").Value = 2688.019336502487
$ws.Range("G6").Value = 18513.97190618839
$ws.Range("H6").Value = 32778.57395031787
$ws.Range("I6").Value = 92217.43862775799
$ws.Range("L6").Value = 4823800.772312711
$ws.Range("M6").Value = 18894452
$ws.Range("N6").Value = -83150891.75703137

# Row 7
$ws.Range("B7").Value = 176563.266501961
$ws.Range("C7").Value = 572392.7046765155
$ws.Range("D7").Value = 16320.09426220185
$ws.Range("E7").Value = 253492.7913824112
$ws.Range("H7").Value = 3153.720010741634
$ws.Range("I7").Value = 28673.69373722013
$ws.Range("L7").Value = 430737.9390276206
$ws.Range("M7").Value = 1293199
$ws.Range("N7").Value = -5090446.900135092

# Row 8
$ws.Range("B8").Value = 73921.97731204861
$ws.Range("C8").Value = 868547.6463691593
$ws.Range("D8").Value = 13332.5334788601
$ws.Range("E8").Value = 323535.1734531521
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 41934.36343044306
$ws.Range("L8").Value = 553276.668324528
$ws.Range("M8").Value = 1750829
$ws.Range("N8").Value = 596773.9453163131

# Row 9
$ws.Range("B9").Value = 127835.8431337425
$ws.Range("C9").Value = 1159158.538727236
$ws.Range("D9").Value = 31289.61305831865
$ws.Range("E9").Value = 426125.1041106441
$ws.Range("H9").Value = 5229.056286106441
$ws.Range("I9").Value = 47618.0998232898
$ws.Range("K9").Value = 98997.98822051805
$ws.Range("L9").Value = 751842.3456269886
$ws.Range("M9").Value = 2128893
$ws.Range("N9").Value = 1699719.576018669

# --- Delete row 10 entirely (was state index 9 / "NC") ---
$ws.Rows.Item(10).Delete()

# --- Update the dimension to reflect the new used range A1:P9 ---
$ws.Range("A1:P9").Select()
